$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data row is inserted at row 320 (pushing the existing rows 320-399
# down to 321-400). The new row carries a fresh record for the same
# market/product (Feria Lagunitas de Puerto Montt - Ciboulette).
$ws.Rows.Item(320).Insert()

$row = $ws.Rows.Item(320)
$row.Cells.Item(1, 1).Value = 4
$row.Cells.Item(1, 2).Value = "Feria Lagunitas de Puerto Montt"
$row.Cells.Item(1, 3).Value = "Los Lagos"
$row.Cells.Item(1, 4).Value = 45173
$row.Cells.Item(1, 5).Value = 10
$row.Cells.Item(1, 6).Value = 100112039
$row.Cells.Item(1, 7).Value = "Ciboulette"
$row.Cells.Item(1, 8).Value = "Sin especificar"
$row.Cells.Item(1, 9).Value = "Primera"
$row.Cells.Item(1, 10).Value = 80
$row.Cells.Item(1, 11).Value = 3500
$row.Cells.Item(1, 12).Value = 3500
$row.Cells.Item(1, 13).Value = 3500
$row.Cells.Item(1, 14).Value = "`$/docena de atados"
$row.Cells.Item(1, 15).Value = "Región Metropolitana"
$row.Cells.Item(1, 16).Value = 1167
$row.Cells.Item(1, 17).Value = 3
$row.Cells.Item(1, 18).Value = "Hortaliza"
